$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-13, columns A-T (matches updated TPM values and
# the reshuffled "ECs" cluster that is now inserted as a new Sending cluster
# group, in addition to updated values for the pre-existing FAPs/MuSCs groups).
$data = @(
    @("ECs", "Il7", "Il2rg", "ECs", 3, 1, 0.8062, 2.4186, 0.581704519790725, 0.5817045197907249, 3, 1, 35.32109533333334, 105.963286, 0.4123245624288747, 0.4123245624288747, 28.47586705773334, 256.2828035196, 0.2398510615856094, 0.2398510615856093),
    @("ECs", "Il7", "Il2rg", "FAPs", 3, 1, 0.8062, 2.4186, 0.581704519790725, 0.5817045197907249, 3, 1, 1.228158333333333, 3.684475, 0.01433703690686912, 0.01433703690686912, 0.9901412483333333, 8.911271235000001, 0.008339919169132202, 0.008339919169132201),
    @("ECs", "Il7", "Il2rg", "MuSCs", 3, 1, 0.8062, 2.4186, 0.581704519790725, 0.5817045197907249, 3, 1, 2.583168333333334, 7.749505, 0.03015489023401347, 0.03015489023401347, 2.082550310333334, 18.742952793, 0.01754123594291883, 0.01754123594291883),
    @("ECs", "Il7", "Il2rg", "Resolving-Mac", 3, 1, 0.8062, 2.4186, 0.581704519790725, 0.5817045197907249, 3, 1, 46.53090866666667, 139.592726, 0.5431835104302428, 0.5431835104302427, 37.51321856706667, 337.6189671036, 0.3159723030930647, 0.3159723030930646),
    @("FAPs", "Il7", "Il2rg", "ECs", 2, 0.6666666666666666, 0.4454536666666667, 1.336361, 0.3214120705251191, 0.3214120705251191, 3, 1, 35.32109533333334, 105.963286, 0.4123245624288747, 0.4123245624288747, 15.73391142691623, 141.605202842246, 0.1325260913386284, 0.1325260913386283),
    @("FAPs", "Il7", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.4454536666666667, 1.336361, 0.3214120705251191, 0.3214120705251191, 3, 1, 1.228158333333333, 3.684475, 0.01433703690686912, 0.01433703690686912, 0.5470876328305555, 4.923788695475, 0.004608096717431853, 0.004608096717431853),
    @("FAPs", "Il7", "Il2rg", "MuSCs", 2, 0.6666666666666666, 0.4454536666666667, 1.336361, 0.3214120705251191, 0.3214120705251191, 3, 1, 2.583168333333334, 7.749505, 0.03015489023401347, 0.03015489023401347, 1.150681805700556, 10.356136251305, 0.009692145706571963, 0.009692145706571962),
    @("FAPs", "Il7", "Il2rg", "Resolving-Mac", 2, 0.6666666666666666, 0.4454536666666667, 1.336361, 0.3214120705251191, 0.3214120705251191, 3, 1, 46.53090866666667, 139.592726, 0.5431835104302428, 0.5431835104302427, 20.72736387889845, 186.546274910086, 0.174585736762487, 0.174585736762487),
    @("MuSCs", "Il7", "Il2rg", "ECs", 2, 0.6666666666666666, 0.1342733333333333, 0.40282, 0.09688340968415604, 0.09688340968415604, 3, 1, 35.32109533333334, 105.963286, 0.4123245624288747, 0.4123245624288747, 4.742681207391112, 42.68413086652, 0.03994740950463704, 0.03994740950463704),
    @("MuSCs", "Il7", "Il2rg", "FAPs", 2, 0.6666666666666666, 0.1342733333333333, 0.40282, 0.09688340968415604, 0.09688340968415604, 3, 1, 1.228158333333333, 3.684475, 0.01433703690686912, 0.01433703690686912, 0.1649089132777778, 1.4841802195, 0.001389021020305066, 0.001389021020305066),
    @("MuSCs", "Il7", "Il2rg", "MuSCs", 2, 0.6666666666666666, 0.1342733333333333, 0.40282, 0.09688340968415604, 0.09688340968415604, 3, 1, 2.583168333333334, 7.749505, 0.03015489023401347, 0.03015489023401347, 0.3468506226777778, 3.1216556041, 0.002921508584522683, 0.002921508584522683),
    @("MuSCs", "Il7", "Il2rg", "Resolving-Mac", 2, 0.6666666666666666, 0.1342733333333333, 0.40282, 0.09688340968415604, 0.09688340968415604, 3, 1, 46.53090866666667, 139.592726, 0.5431835104302428, 0.5431835104302427, 6.247860209702222, 56.23074188731999, 0.05262547057469126, 0.05262547057469125)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}
